# Add a new product ("truskawka" / strawberry) with its calorie value (5)
# to the two linked sheets, then leave the UI selection/active-sheet state
# the way the author left it when they saved: cursor parked on the new
# row of "produkty", and "kalorie" as the active/visible tab with the
# cursor resting on B8.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("produkty")
$ws2 = $wb.Worksheets.Item("kalorie")

# New data row on both sheets (row 6).
$ws1.Range("A6").Value = "truskawka"
$ws2.Range("A6").Value = 5

# Selection / active-sheet bookkeeping, applied in the order that
# produces the saved end state: select A6 on "produkty" first, then
# switch to "kalorie" and select B8 there, so "kalorie" ends up as the
# active tab.
$ws1.Range("A6").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("B8").Select() | Out-Null
